# Remove "Francisco Vidal" (row 16) and "Jorge Tarud" (row 18) from the
# candidates table (both "Unidad Constituyente" / PPD entries). Deleting
# row 18 first keeps row 16's index valid for the second delete.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Delete()
$ws.Rows.Item(16).Delete()

# Column A ("id") is a static, pre-filled 1..N sequence rather than a
# formula, so restore it to a plain consecutive run after the rows above
# shifted up.
For ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$ws.Range("A2:A18").Select()
